# Weekly market-data refresh for the Leve profit sheets.
# Updates cached NQ/HQ price + profit columns (H:N) on ALC, ARM, BSM, CRP,
# CUL with the latest Universalis snapshot; some leves did not have any
# cached price data yet and get it filled in for the first time here.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")  # row 74
$ws.Range("H74").Value = 8991.143
$ws.Range("I74").Value = 8994.5
$ws.Range("J74").Value = 8989.799999999999
$ws.Range("K74").Value = 8994.5
$ws.Range("L74").Value = 8989.799999999999
$ws.Range("M74").Value = -8058.5
$ws.Range("N74").Value = -10861.8

$ws = $wb.Worksheets.Item("ALC")  # row 77
$ws.Range("H77").Value = 8991.143
$ws.Range("I77").Value = 8994.5
$ws.Range("J77").Value = 8989.799999999999
$ws.Range("K77").Value = 44972.5
$ws.Range("L77").Value = 44949
$ws.Range("M77").Value = -40292.5
$ws.Range("N77").Value = -54309

$ws = $wb.Worksheets.Item("ALC")  # row 92
$ws.Range("H92").Value = 587.4
$ws.Range("I92").Value = 587.4
$ws.Range("K92").Value = 587.4
$ws.Range("M92").Value = 660.6

$ws = $wb.Worksheets.Item("ALC")  # row 112
$ws.Range("H112").Value = 6088.8
$ws.Range("J112").Value = 7361
$ws.Range("L112").Value = 22083
$ws.Range("N112").Value = -24299

$ws = $wb.Worksheets.Item("ALC")  # row 125
$ws.Range("H125").Value = 2900
$ws.Range("I125").Value = 2900
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 26100
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -23640

$ws = $wb.Worksheets.Item("ALC")  # row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0

$ws = $wb.Worksheets.Item("ALC")  # row 127
$ws.Range("H127").Value = 1107.5
$ws.Range("I127").Value = 715
$ws.Range("J127").Value = 1500
$ws.Range("K127").Value = 2145
$ws.Range("L127").Value = 4500
$ws.Range("M127").Value = 2815
$ws.Range("N127").Value = -14420

$ws = $wb.Worksheets.Item("ALC")  # row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0

$ws = $wb.Worksheets.Item("ALC")  # row 129
$ws.Range("H129").Value = 4000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 4000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 12000
$ws.Range("N129").Value = -22000

$ws = $wb.Worksheets.Item("ALC")  # row 130
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0

$ws = $wb.Worksheets.Item("ALC")  # row 131
$ws.Range("H131").Value = 8224.75
$ws.Range("I131").Value = 7299.6665
$ws.Range("J131").Value = 11000
$ws.Range("K131").Value = 21898.9995
$ws.Range("L131").Value = 33000
$ws.Range("M131").Value = -16858.9995
$ws.Range("N131").Value = -43080

$ws = $wb.Worksheets.Item("ALC")  # row 132
$ws.Range("H132").Value = 4014
$ws.Range("I132").Value = 4017.111
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 12051.333
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -9521.332999999999
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("ALC")  # row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0

$ws = $wb.Worksheets.Item("ALC")  # row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0

$ws = $wb.Worksheets.Item("ALC")  # row 135
$ws.Range("H135").Value = 966
$ws.Range("I135").Value = 966
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8694
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -6159

$ws = $wb.Worksheets.Item("ALC")  # row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0

$ws = $wb.Worksheets.Item("ALC")  # row 137
$ws.Range("H137").Value = 6799.6
$ws.Range("I137").Value = 6799.6
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 20398.8
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -17848.8

$ws = $wb.Worksheets.Item("ALC")  # row 138
$ws.Range("H138").Value = 3861.7778
$ws.Range("I138").Value = 1000
$ws.Range("J138").Value = 4679.4287
$ws.Range("K138").Value = 3000
$ws.Range("L138").Value = 14038.2861
$ws.Range("M138").Value = 2140
$ws.Range("N138").Value = -24318.2861

$ws = $wb.Worksheets.Item("ALC")  # row 139
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0

$ws = $wb.Worksheets.Item("ALC")  # row 140
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0

$ws = $wb.Worksheets.Item("ALC")  # row 141
$ws.Range("H141").Value = 1500
$ws.Range("I141").Value = 2000
$ws.Range("J141").Value = 1000
$ws.Range("K141").Value = 6000
$ws.Range("L141").Value = 3000
$ws.Range("M141").Value = -820
$ws.Range("N141").Value = -13360

$ws = $wb.Worksheets.Item("ARM")  # row 121
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0

$ws = $wb.Worksheets.Item("ARM")  # row 122
$ws.Range("H122").Value = 2501.625
$ws.Range("I122").Value = 1499.75
$ws.Range("J122").Value = 3503.5
$ws.Range("K122").Value = 4499.25
$ws.Range("L122").Value = 10510.5
$ws.Range("M122").Value = -2049.25
$ws.Range("N122").Value = -15410.5

$ws = $wb.Worksheets.Item("ARM")  # row 123
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0

$ws = $wb.Worksheets.Item("ARM")  # row 124
$ws.Range("H124").Value = 25000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 25000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 25000
$ws.Range("N124").Value = -34820

$ws = $wb.Worksheets.Item("ARM")  # row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0

$ws = $wb.Worksheets.Item("ARM")  # row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0

$ws = $wb.Worksheets.Item("ARM")  # row 127
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0

$ws = $wb.Worksheets.Item("ARM")  # row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0

$ws = $wb.Worksheets.Item("ARM")  # row 129
$ws.Range("H129").Value = 18899
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 18899
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 18899
$ws.Range("N129").Value = -28899

$ws = $wb.Worksheets.Item("ARM")  # row 130
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0

$ws = $wb.Worksheets.Item("ARM")  # row 131
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0

$ws = $wb.Worksheets.Item("ARM")  # row 132
$ws.Range("H132").Value = 2625.3044
$ws.Range("I132").Value = 1105.4615
$ws.Range("J132").Value = 4601.1
$ws.Range("K132").Value = 3316.3845
$ws.Range("L132").Value = 13803.3
$ws.Range("M132").Value = -786.3844999999997
$ws.Range("N132").Value = -18863.3

$ws = $wb.Worksheets.Item("ARM")  # row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0

$ws = $wb.Worksheets.Item("ARM")  # row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0

$ws = $wb.Worksheets.Item("ARM")  # row 135
$ws.Range("H135").Value = 100000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 100000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140

$ws = $wb.Worksheets.Item("ARM")  # row 137
$ws.Range("H137").Value = 50000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 50000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws = $wb.Worksheets.Item("ARM")  # row 138
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0

$ws = $wb.Worksheets.Item("ARM")  # row 139
$ws.Range("H139").Value = 119997.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 119997.5
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 119997.5
$ws.Range("N139").Value = -130277.5

$ws = $wb.Worksheets.Item("ARM")  # row 140
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0

$ws = $wb.Worksheets.Item("ARM")  # row 141
$ws.Range("H141").Value = 40000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 40000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360

$ws = $wb.Worksheets.Item("BSM")  # row 130
$ws.Range("H130").Value = 100000
$ws.Range("J130").Value = 100000
$ws.Range("L130").Value = 100000
$ws.Range("N130").Value = -110040

$ws = $wb.Worksheets.Item("CRP")  # row 68
$ws.Range("H68").Value = 30000
$ws.Range("I68").Value = 30000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 30000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -29251
$ws.Range("N68").ClearContents()

$ws = $wb.Worksheets.Item("CRP")  # row 71
$ws.Range("H71").Value = 30000
$ws.Range("I71").Value = 30000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 90000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -86256
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("CRP")  # row 74
$ws.Range("H74").Value = 63316.668
$ws.Range("J74").Value = 63316.668
$ws.Range("L74").Value = 63316.668
$ws.Range("N74").Value = -65064.668

$ws = $wb.Worksheets.Item("CRP")  # row 77
$ws.Range("H77").Value = 63316.668
$ws.Range("J77").Value = 63316.668
$ws.Range("L77").Value = 189950.004
$ws.Range("N77").Value = -198686.004

$ws = $wb.Worksheets.Item("CUL")  # row 104
$ws.Range("H104").Value = 4999.864
$ws.Range("J104").Value = 4999.864
$ws.Range("L104").Value = 14999.592
$ws.Range("N104").Value = -20241.592

$ws = $wb.Worksheets.Item("CUL")  # row 121
$ws.Range("H121").Value = 1148.75
$ws.Range("J121").Value = 1148.75
$ws.Range("L121").Value = 3446.25
$ws.Range("N121").Value = -6066.25

